$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Oct 26 2020 vs Kings XI Punjab) is replaced with the data that used
# to live in row 4 (Oct 7 2020 @ Abu Dhabi vs Chennai Super Kings); rows 3-5
# (old rows 3, 4, 5) are then deleted outright, shrinking the sheet to just
# the header row plus this single data row.
#
# D2 ("Kolkata Knight Riders"), F2 (player name) and I2/J2 (both "0") already
# hold the same text that the surviving row needs, so they are left
# untouched to avoid any unintended re-encoding (F2 in particular carries a
# trailing NBSP, not a plain space).

# Force text storage for the numeric-looking columns so they stay as strings
# (matching the original "numberStoredAsText" cells) instead of being
# auto-converted to real numbers. Style is restored to "Normal" afterwards so
# no stray number-format style is left behind on the cells.
$numCols = $ws.Range("G2:K2")
$numCols.NumberFormat = "@"

$ws.Range("A2").Value = " Oct 7 2020"
$ws.Range("B2").Value = " Abu Dhabi"
$ws.Range("C2").Value = "KKR won by 10 runs"
$ws.Range("E2").Value = "Chennai Super Kings"
$ws.Range("G2").Value = "0"
$ws.Range("H2").Value = "2"
$ws.Range("K2").Value = "0.00"

$numCols.Style = "Normal"

# Remove the now-obsolete rows 3-5, shrinking the used range down to A1:K2.
$ws.Rows("3:5").Delete()
